# Update the sampling-method lookup table (dbo_sampmet):
# remove the erroneous "?" / "DCN" row (sampmet_id = 22) and let the
# following rows shift up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 holds sampmet_id=22, sampmet="?", code="DCN" -> delete it entirely.
$ws.Rows.Item(24).Delete() | Out-Null

# Keep the named range "dbo_sampmet" in sync with the now-smaller table
# (it covered A1:G27, now the last row is 26).
$wb.Names.Item("dbo_sampmet").RefersTo = "=dbo_sampmet!`$A`$1:`$G`$26"

# Restore the selected cell as left by the editor.
$ws.Range("B17").Select() | Out-Null
